$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I = "APISentiment". Rows whose value becomes 0
$rowsToZero = @(
    2,4,8,9,10,11,12,13,15,16,19,22,23,24,25,28,32,33,34,35,
    36,37,39,40,43,46,47,48,49,50,51,52,56,57,58,59,60,61,63,64,
    67,70,71,72,73,74,80,81,82,83,84,85,87,88,91,94,95,96,97,98,
    104,105,106,107,108,109,111,112,115,118,119,120,121,122,128,129,130,131,132,133,
    135,136,139,142,143,144,145,146,150,151,152,153,154,155,156,157,159,160,163,166,
    167,168,169,170,174,175,176,177,178,179,180,181,183,184,187,190,191,192,193,194,
    198,199,201,202,203,204,205,207,208,211,214,215,216,217,218,222,223,226,227,228,
    229,231,232,235,238,239,240,241,242,246,247,251,252,253,255,256,259,262,263,264,
    265,266,270,271,275,276,277,279,280,283,286,287,288,289,290,294,295,299,300,303,
    304,307,310,311,312,313,314,318,319,323,324,326,327,328,331,334,335,336,337,338,
    342,343,347,348,351,352,355,358,359,360,361,362,366,367,371,372,375,379,382,383,
    384,385,386,390,391,395,396,399,400,403,406,407,408,409,410,414,415,419,420,423,
    424,426,427,430,431,432,433,434,438,439,443,444,447,448,454,455,456,457,458,462,
    463,467,468,471,472,475,478,479,480,481,482,486,487,491,492,495,496,499,500,501,
    502,503,504,505,506,510,511,515,516,519,520,523,524,526,527,528,529,530,534,535,
    539,540,543,544,547,548,551,552,553,554,558,559,563,564,567,568,571,572
)

# Rows whose value becomes 0.995 (previously 0.997)
$rowsTo995 = @(18,42,66,90,114,138,162,186,210,234,258,282,306,330,354,378,401,402)

foreach ($r in $rowsToZero) {
    $ws.Cells.Item($r, 9).Value = 0
}

foreach ($r in $rowsTo995) {
    $ws.Cells.Item($r, 9).Value = 0.995
}
